$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 12

# Column A holds a date-formatted-looking string ("2025-09-08"). Assigning it
# directly would make Excel auto-detect it as a real date and reformat the
# cell, which the source data does not want (it must stay a plain text
# value, same as every other cell in this sheet). Force text interpretation
# by applying a text number format before assigning the value, then drop the
# number format again so the cell keeps the sheet's default (unstyled) look.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value2 = "2025-09-08"
$cellA.ClearFormats()

$ws.Cells.Item($row, 2).Value = "21:21:01"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1620.1186 ARS"
